$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 did not previously exist in the sheet, so first bring over the same
# cell formatting used by the data rows above (row 4) before writing values,
# matching what Excel does when a user fills in a previously-blank row.
$ws.Range("A4:J4").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 already existed (blank placeholder row) so its formatting is already
# in place; just populate the values.
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Select the newly-entered rows, mirroring what the author's Excel session
# left selected when it saved the file.
$ws.Range("A5:J6").Select() | Out-Null
